$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended below the existing manufacturer/location table
# (Manufacturer, Material, Latitude, Longitude)
$ws.Range("A6").Value = "GSM Industrial"
$ws.Range("B6").Value = "Steel"

# Latitude: plain text, same "numeric-looking string stored as text" pattern
# used by the rest of the Latitude/Longitude columns.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "53.493746"
$ws.Range("C6").ClearFormats()

# Longitude: also stored as text, but this one keeps an explicit "General"
# number format applied to it (picks up a new style entry).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "-2.2086889"
$ws.Range("D6").NumberFormat = "General"

# Match the trailing UI state left behind in the saved file.
$ws.Range("I8").Select() | Out-Null
